# Correct most names to the official names from the website.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G6").Value = "Ramanagara"
$ws.Range("G13").Value = "Vijayapura (Bijapur)"
$ws.Range("F19").Value = ""
$ws.Range("G20").Value = "Davangere"
$ws.Range("G21").Value = "Davangere"
$ws.Range("G25").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G28").Value = "Davangere"
$ws.Range("G29").Value = "Davangere"
$ws.Range("G31").Value = "Vijayapura (Bijapur)"
$ws.Range("F37").Value = ""
$ws.Range("G40").Value = "Davangere"
$ws.Range("G41").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G45").Value = "Vijayapura (Bijapur)"
$ws.Range("G46").Value = "Vijayapura (Bijapur)"
$ws.Range("G51").Value = "Kalaburagi (Gulbarga)"
$ws.Range("G57").Value = "Vijayapura (Bijapur)"
$ws.Range("G58").Value = "Kalaburagi (Gulbarga)"
